$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 ("Operating revenue (Turnover)") - this shifts all
# subsequent rows up by one and removes the now-unused shared string.
$ws.Rows.Item(3).Delete()

# Move the selection as recorded after the edit.
$ws.Range("C17").Select()
